$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row above the current row 67 ("FoOMCtiL"), shifting it and
# everything below down by one. Excel copies the formatting of the row
# above (row 66) into the freshly inserted row, matching native behavior.
$ws.Rows.Item(67).Insert()

# Populate the new row with the "ETCCwTC" acronym entry (Elasticity of
# Transmission Connectivity Coefficient wrt Transmission Capacity), under
# the "elec" top level folder, with "medium" importance to update.
$ws.Cells.Item(67, 1).Value = "elec"
$ws.Cells.Item(67, 2).Value = "ETCCwTC"
$ws.Cells.Item(67, 3).Value = "Elasticity of Transmission Connectivity Coefficient wrt Transmission Capacity"
$ws.Cells.Item(67, 6).Value = "medium"

# Match the "medium" importance cell's fill/format used elsewhere in the
# column (row 68, the shifted-down "FoOMCtiL" row, which already carries
# the correct "medium" formatting) rather than the inherited format from
# row 66 ("optional").
$ws.Cells.Item(68, 6).Copy()
$ws.Cells.Item(67, 6).PasteSpecial(-4122)

# Row 66 had a value in column G; the inserted row inherited that
# (empty) cell's formatting. Remove it since the new row has no entry
# for that column.
$ws.Cells.Item(67, 7).Clear()

# The row's text wraps onto two lines, so give it the corresponding
# row height (matches the other two-line rows in this sheet).
$ws.Rows.Item(67).RowHeight = 29
